$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (F column) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 257
$ws1.Range("F5").Value = 1686
$ws1.Range("F6").Value = 1486
$ws1.Range("F8").Value = 60
$ws1.Range("F9").Value = 438
$ws1.Range("F10").Value = 118

# Sheet "演出" (Performances) - update "想去人数" (F column) count
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 19

# Sheet "全部类型" (All types) - update "想去人数" (F column) counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 257
$ws4.Range("F5").Value = 1686
$ws4.Range("F6").Value = 1486
$ws4.Range("F8").Value = 19
$ws4.Range("F9").Value = 60
$ws4.Range("F10").Value = 438
$ws4.Range("F11").Value = 118
